# Fruta / hortaliza, semanal
# Shuffle the weekly price rows (2-8) into their new order and append a
# new trailing row (9) that holds the data that used to be in row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for each data row: Fecha (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M), Precio $/Kg (P)
$rows = @(
    @{ Row = 2; Fecha = 44243; Volumen = 1200; Min = 1200; Max = 1300; Prom = 1250; PrecioKg = 1250 },
    @{ Row = 3; Fecha = 44229; Volumen = 1500; Min = 1400; Max = 1500; Prom = 1450; PrecioKg = 1450 },
    @{ Row = 4; Fecha = 44407; Volumen = 1000; Min = 1200; Max = 1300; Prom = 1250; PrecioKg = 1250 },
    @{ Row = 5; Fecha = 44284; Volumen = 1500; Min =  800; Max =  850; Prom =  825; PrecioKg =  825 },
    @{ Row = 6; Fecha = 44291; Volumen = 1000; Min = 1000; Max = 1200; Prom = 1100; PrecioKg = 1100 },
    @{ Row = 7; Fecha = 44341; Volumen = 1300; Min =  900; Max = 1000; Prom =  950; PrecioKg =  950 },
    @{ Row = 8; Fecha = 44442; Volumen = 1250; Min =  850; Max =  900; Prom =  875; PrecioKg =  875 },
    @{ Row = 9; Fecha = 44175; Volumen = 1600; Min = 1000; Max = 1200; Prom = 1100; PrecioKg = 1100 }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("D$n").Value = $r.Fecha
    $ws.Range("J$n").Value = $r.Volumen
    $ws.Range("K$n").Value = $r.Min
    $ws.Range("L$n").Value = $r.Max
    $ws.Range("M$n").Value = $r.Prom
    $ws.Range("P$n").Value = $r.PrecioKg
}

# Row 9 is brand new - fill in the columns that stay constant across all
# rows for this sheet (copied from row 6, which used to hold this data),
# and make sure the date cell picks up the same date number format.
$ws.Range("A9").Value = $ws.Range("A6").Value()
$ws.Range("B9").Value = $ws.Range("B6").Value()
$ws.Range("C9").Value = $ws.Range("C6").Value()
$ws.Range("E9").Value = $ws.Range("E6").Value()
$ws.Range("F9").Value = $ws.Range("F6").Value()
$ws.Range("G9").Value = $ws.Range("G6").Value()
$ws.Range("H9").Value = $ws.Range("H6").Value()
$ws.Range("I9").Value = $ws.Range("I6").Value()
$ws.Range("N9").Value = $ws.Range("N6").Value()
$ws.Range("O9").Value = $ws.Range("O6").Value()
$ws.Range("Q9").Value = $ws.Range("Q6").Value()
$ws.Range("R9").Value = $ws.Range("R6").Value()

$ws.Range("D9").NumberFormat = $ws.Range("D2").NumberFormat
